# chore: update Sheets via scheduled runner
# Applies refreshed currentAveragePrice / LevePrice / LeveProfit figures
# (columns H, I, J, K, L, M, N) to a handful of leve rows across several
# job sheets, as produced by the scheduled market-price refresh.

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        [double]$H,
        [double]$I,
        [double]$J,
        [double]$K,
        [double]$L,
        [double]$M,
        [double]$N,
        [bool]$SetM
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Cells.Item($Row, 8).Value  = $H   # H - currentAveragePrice
    $ws.Cells.Item($Row, 9).Value  = $I   # I - currentAveragePriceNQ
    $ws.Cells.Item($Row, 10).Value = $J   # J - currentAveragePriceHQ
    $ws.Cells.Item($Row, 11).Value = $K   # K - LevePriceNQ
    $ws.Cells.Item($Row, 12).Value = $L   # L - LevePriceHQ
    if ($SetM) {
        $ws.Cells.Item($Row, 13).Value = $M   # M - LeveProfitNQ
    }
    $ws.Cells.Item($Row, 14).Value = $N   # N - LeveProfitHQ
}

# ---- ALC ----
Set-LeveRow "ALC" 116 4596.6665 5590.8 3886.5715 5590.8 3886.5715 -2148.8 -10770.5715 $true
Set-LeveRow "ALC" 129 984.88 509 1135.1578 1527 3405.4734 3473 -13405.4734 $true
Set-LeveRow "ALC" 138 2045219.2 3671 3973348.2 11013 11920044.6 -5873 -11930324.6 $true

# ---- ARM ----
Set-LeveRow "ARM" 32  8478.413 4520.2114 27189.908 4520.2114 27189.908 -4233.2114 -27763.908 $true
Set-LeveRow "ARM" 61  2426.6 2068.75 3858 2068.75 3858 -1856.75 -4282 $true
Set-LeveRow "ARM" 122 1952.5306 1715.3667 2327 5146.1001 6981 -2696.1001 -11881 $true
Set-LeveRow "ARM" 132 1659.909 1408.3529 2515.2 4225.0587 7545.599999999999 -1695.0587 -12605.6 $true
Set-LeveRow "ARM" 136 2426.6 2068.75 3858 6206.25 11574 -3656.25 -16674 $true

# ---- CRP ----
Set-LeveRow "CRP" 31  1813.4722 1336.65 2409.5 1336.65 2409.5 -1041.65 -2999.5 $true
Set-LeveRow "CRP" 34  1813.4722 1336.65 2409.5 1336.65 2409.5 -1134.65 -2813.5 $true
Set-LeveRow "CRP" 58  2427.6538 1308.7916 3386.6785 1308.7916 3386.6785 -1105.7916 -3792.6785 $true
Set-LeveRow "CRP" 132 5540.7334 6414.25 4542.4287 19242.75 13627.2861 -16712.75 -18687.2861 $true
Set-LeveRow "CRP" 136 2427.6538 1308.7916 3386.6785 3926.3748 10160.0355 -1376.3748 -15260.0355 $true

# ---- CUL ----
Set-LeveRow "CUL" 68  945.7808 698.2059 1161.6154 2094.6177 3484.8462 -1283.6177 -5106.8462 $true
Set-LeveRow "CUL" 71  945.7808 698.2059 1161.6154 6283.8531 10454.5386 -2227.8531 -18566.5386 $true
Set-LeveRow "CUL" 107 24425.191 18657.582 35363.758 55972.746 106091.274 -54052.746 -109931.274 $true
Set-LeveRow "CUL" 114 1487.2174 690.375 3308.5715 2071.125 9925.7145 1182.875 -16433.7145 $true
Set-LeveRow "CUL" 137 34305.574 2787.9333 60570.277 8363.7999 181710.831 -3263.7999 -191910.831 $true

# ---- GSM ----
Set-LeveRow "GSM" 80  65224.316 2086 135378 2086 135378 -1088 -137374 $true
Set-LeveRow "GSM" 83  65224.316 2086 135378 10430 676890 -5438 -686874 $true
Set-LeveRow "GSM" 132 1845.8937 2001.3572 1779.9395 6004.071599999999 5339.818499999999 -3474.071599999999 -10399.8185 $true

# ---- LTW ----
Set-LeveRow "LTW" 16 1951.65 2031.2632 439 2031.2632 439 -1861.2632 -779 $true
Set-LeveRow "LTW" 22 293.46667 280 320.4 280 320.4 15 -910.4 $true
Set-LeveRow "LTW" 27 293.46667 280 320.4 280 320.4 -173 -534.4 $true
Set-LeveRow "LTW" 46 1010.4 1600 757.7143 1600 757.7143 -1412 -1133.7143 $true

# ---- WVR (LeveProfitNQ / column M not present on this row; leave untouched) ----
Set-LeveRow "WVR" 133 41000 0 41000 0 41000 0 -51120 $false
